$wb = $excel.ActiveWorkbook

$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (this single shared string is referenced by Overview!E2/F2/E3/F3 and by
#  the "Status" column on both language sheets, so updating any one cell
#  to the new text and letting Excel's de-dup reuse the string updates
#  them all; we set each cell explicitly to be safe.)
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

# --- Latest Handback DateTime refreshed, Error Detail cleared (handback now in sync) ---
$zh.Range("K2").Value = "2016-08-31 13:09:57"
$zh.Range("K3").Value = "2016-08-31 13:09:57"
$zh.Range("P2").Value = ""

$de.Range("K2").Value = "2016-08-31 13:10:35"
$de.Range("K3").Value = "2016-08-31 13:10:35"
$de.Range("P2").Value = ""

# --- Column width adjustments (report regenerated with wider Status / narrower Error Detail columns) ---
# (ColumnWidth is quantized by the host to 1/6-character steps on write, so these
#  inputs are chosen to land the stored OOXML width on the nearest achievable step
#  to the target widths 29.9777050018311 / 13.7470531463623.)
$ov.Range("E1").ColumnWidth = 29.15
$ov.Range("F1").ColumnWidth = 29.15

$zh.Range("C1").ColumnWidth = 29.15
$zh.Range("P1").ColumnWidth = 12.85

$de.Range("C1").ColumnWidth = 29.15
$de.Range("P1").ColumnWidth = 12.85
